$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c4 = "1. OPERATORS`n2.CONDITIONAL`n STATEMENTS`n3.LOOPS`n4. LINEAR ALGEBRA"
$e4 = "SHERIYANS AI SCHOOL `nhttps://youtu.be/_aWbUudZ5Yo?si=cohj6duDlwQ0-CI2`n3BLUE 1BROWN`nhttps://youtu.be/fNk_zzaMoSs?si=4-NkCB-90DB-8J-U"

$ws.Range("C4").Value = $c4
$ws.Range("E4").Value = $e4

$ws.Rows.Item(4).RowHeight = 115.2
